$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.919.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.284.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.69%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +5.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0962"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.56%  "

$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.623.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.820"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.300.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.794.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("E29").Value = "  -1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.82%  "

$ws.Range("E32").Value = "  +4.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0700"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("E40").Value = "  +3.52%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +32.64%  "

$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000225"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.488.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.75%  "
